$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.582.66"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.287.65"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("E7").Value = "  -3.36%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -3.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0783"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.25%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").Value = "2.643.22"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "2.284.00"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "42.491.17"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "0.0₃0889"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E24").Value = "  -0.51%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0684"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("E42").Value = "  -3.37%  "
$ws.Range("D43").Value = "1.992.20"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.02%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.42%  "
$ws.Range("E47").Value = "  -10.10%  "
$ws.Range("E48").Value = "  -1.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("D51").Value = "2.510.33"
$ws.Range("E51").Value = "  -0.39%  "
